$wb = $excel.ActiveWorkbook

# 1) Create the new "Creative_Conversion_S3_Mapper" sheet by copying the
#    structurally-closest existing sheet (Creative_Delivery_S3_Mapper) and
#    placing it right after "Creative_Conversion_Mapper".
$template = $wb.Worksheets.Item("Creative_Delivery_S3_Mapper")
$anchor = $wb.Worksheets.Item("Creative_Conversion_Mapper")
$template.Copy($null, $anchor)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "Creative_Conversion_S3_Mapper"

# 2) The new sheet should describe creative CONVERSION data (27 rows), not
#    delivery data (28 rows) - drop the trailing "Cost" row and swap the
#    last two metric rows for the conversion-specific columns.
$newSheet.Rows.Item(28).Delete()
$newSheet.Range("A26").Value = "Click_Based_Conversions"
$newSheet.Range("B26").Value = "Click_Based_Conversions"
$newSheet.Range("A27").Value = "Impression_Based_Conversions"
$newSheet.Range("B27").Value = "Impression_Based_Conversions"

# 3) Fix up the per-sheet selections / active-cell state to match the edit.
$convMapper = $wb.Worksheets.Item("Creative_Conversion_Mapper")
$convMapper.Activate()
$convMapper.Range("A1:E27").Select()

$deliveryS3 = $wb.Worksheets.Item("Creative_Delivery_S3_Mapper")
$deliveryS3.Activate()
$deliveryS3.Range("G16").Select()

# 4) Leave the new sheet as the active / selected tab.
$newSheet.Activate()
$newSheet.Range("F9").Select()
